# Fix the numbering of the uncertainty table on the Unified_table sheet.
# Row 43 previously duplicated the index number of the preceding block (22);
# it should carry its own unique sequence number (23). Every subsequent
# row's running index is driven by "=A43+1", "=A44+1", ... formulas, so
# they recalculate one higher automatically, all the way down to row 74.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unified_table")

$ws.Range("A43").Value = 23

# The table (and its AutoFilter dropdowns) now cover one more data row
# (1-74 instead of 1-73) - re-apply the filter over the new extent and
# keep the backing _FilterDatabase defined name in sync with it.
$ws.AutoFilterMode = $false
$ws.Range("A1:Q74").AutoFilter()

foreach ($n in $wb.Names) {
  if ($n.Name -eq "Unified_table!_FilterDatabase") {
    $n.RefersTo = "=Unified_table!`$A`$1:`$Q`$74"
  }
}
